$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before I (existing column I "Diễn giải" shifts to J).
# Excel normally copies the left neighbour's column formatting/width on
# insert; this engine copies per-cell styles automatically but not the
# column width, so that is reapplied explicitly below.
$ws.Columns("I").Insert()

# New header cell (row 4) and template placeholder cell (row 6).
$ws.Range("I4").Value = "Hub"
$ws.Range("I6").Value = "{hub}"

# I6 should look like the neighbouring amount columns (copy G6's style)
# but centered instead of left-aligned.
$ws.Range("G6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I6").Locked = $true

# Match column H's width on the freshly inserted column I.
$ws.Columns("I").ColumnWidth = 27.27

# Refresh the autofilter range to include the new column.
$ws.AutoFilterMode = $false
$ws.Range("A5:J5").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the autofilter.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Advance Report!_FilterDatabase") {
        $n.RefersTo = "='Advance Report'!`$A`$5:`$J`$5"
    }
}

# Match the author's final cursor position.
$ws.Range("I6").Select()
